$wb = $excel.ActiveWorkbook

# --- Sheet: y_fitted_on_begin_2016 (column B, rows 2-21) ---
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Range("B2").Value = -0.2938333355280292
$ws1.Range("B3").Value = 15.75091470315856
$ws1.Range("B4").Value = 16.16840410683831
$ws1.Range("B5").Value = 17.24757140726266
$ws1.Range("B6").Value = 16.83135171612134
$ws1.Range("B7").Value = 15.45683179089128
$ws1.Range("B8").Value = 15.40029986381569
$ws1.Range("B9").Value = 14.05631650138167
$ws1.Range("B10").Value = 13.97094909323009
$ws1.Range("B11").Value = 12.81988187335937
$ws1.Range("B12").Value = 12.08232313850581
$ws1.Range("B13").Value = 11.4481561755916
$ws1.Range("B14").Value = 9.902301890852771
$ws1.Range("B15").Value = 9.660569876665109
$ws1.Range("B16").Value = 9.841491634848071
$ws1.Range("B17").Value = 9.754029791885703
$ws1.Range("B18").Value = 9.875648563195281
$ws1.Range("B19").Value = 9.23726498277656
$ws1.Range("B20").Value = 9.187067571689358
$ws1.Range("B21").Value = 10.03263826083776

# --- Sheet: y_pred_on_2017_2021 (columns B,C,D rows 2-6) ---
$ws2 = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws2.Range("B2").Value = 9.73677216079953
$ws2.Range("C2").Value = 8.781617266557998
$ws2.Range("D2").Value = 10.69192705504106

$ws2.Range("B3").Value = 9.485659611132803
$ws2.Range("C3").Value = 7.829595616695726
$ws2.Range("D3").Value = 11.14172360556988

$ws2.Range("B4").Value = 9.20961420282806
$ws2.Range("C4").Value = 6.962783394152994
$ws2.Range("D4").Value = 11.45644501150313

$ws2.Range("B5").Value = 8.923187340484494
$ws2.Range("C5").Value = 6.172169862629597
$ws2.Range("D5").Value = 11.67420481833939

$ws2.Range("B6").Value = 8.632437885632589
$ws2.Range("C6").Value = 5.441839916754424
$ws2.Range("D6").Value = 11.82303585451075

# --- Sheet: y_fitted_on_begin_2021 (column B, rows 2-26) ---
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Range("B2").Value = -0.376768595936489
$ws3.Range("B3").Value = 15.66797847531893
$ws3.Range("B4").Value = 15.97003704622982
$ws3.Range("B5").Value = 16.89164175352523
$ws3.Range("B6").Value = 16.69968260064824
$ws3.Range("B7").Value = 15.56448685604645
$ws3.Range("B8").Value = 15.33982918030426
$ws3.Range("B9").Value = 14.17937660761066
$ws3.Range("B10").Value = 13.92022465174941
$ws3.Range("B11").Value = 12.9109003018646
$ws3.Range("B12").Value = 12.1401887105274
$ws3.Range("B13").Value = 11.47756769447617
$ws3.Range("B14").Value = 10.08831205033546
$ws3.Range("B15").Value = 9.656644118445413
$ws3.Range("B16").Value = 9.705253640691833
$ws3.Range("B17").Value = 9.627457686034848
$ws3.Range("B18").Value = 9.714086311128989
$ws3.Range("B19").Value = 9.202933615498912
$ws3.Range("B20").Value = 9.083710526489201
$ws3.Range("B21").Value = 9.746906637564393
$ws3.Range("B22").Value = 9.603962240260932
$ws3.Range("B23").Value = 9.329154049643474
$ws3.Range("B24").Value = 9.5786772847283
$ws3.Range("B25").Value = 7.058647388845118
$ws3.Range("B26").Value = 6.925918852554233

# --- Sheet: y_pred_on_2022_2026 (columns B,C,D rows 2-6) ---
$ws4 = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws4.Range("B2").Value = 6.515190195251154
$ws4.Range("C2").Value = 5.329707004648185
$ws4.Range("D2").Value = 7.700673385854124

$ws4.Range("B3").Value = 6.137551209427303
$ws4.Range("C3").Value = 4.32185838159057
$ws4.Range("D3").Value = 7.953244037264036

$ws4.Range("B4").Value = 5.760643270052815
$ws4.Range("C4").Value = 3.464427835270868
$ws4.Range("D4").Value = 8.056858704834761

$ws4.Range("B5").Value = 5.383852366189734
$ws4.Range("C5").Value = 2.689011239936245
$ws4.Range("D5").Value = 8.078693492443222

$ws4.Range("B6").Value = 5.00708019890754
$ws4.Range("C6").Value = 1.965052996897389
$ws4.Range("D6").Value = 8.049107400917691
